# update version name to 1.20.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bump the PENGUIN1SA4062 row from 1.18 -> 1.19
$ws.Range("I2").Value = "CQLive-PENGUIN1SA4062-1.19"
$ws.Range("J2").Value = "CQLive-PENGUIN1SA4062-1.19.apk"

# Bump the TCL55A261 row from 1.19 -> 1.20
$ws.Range("I3").Value = "CQLive-TCL55A261-1.20"
$ws.Range("J3").Value = "CQLive-TCL55A261-1.20.apk"

# Move the active selection to J12 (matches the saved view state in the diff)
$null = $ws.Range("J12").Select()
